# Auto-generated Excel COM-interop script applying the OOXML diff
# (scheduled runner update to Ridill_Profits sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3157.9697
$ws.Range("I64").Value = 2925.4285
$ws.Range("J64").Value = 3564.9167
$ws.Range("K64").Value = 2925.4285
$ws.Range("L64").Value = 3564.9167
$ws.Range("M64").Value = -2677.4285
$ws.Range("N64").Value = -4060.9167
$ws.Range("H67").Value = 3157.9697
$ws.Range("I67").Value = 2925.4285
$ws.Range("J67").Value = 3564.9167
$ws.Range("K67").Value = 2925.4285
$ws.Range("L67").Value = 3564.9167
$ws.Range("M67").Value = -2067.4285
$ws.Range("N67").Value = -5280.9167
$ws.Range("H69").Value = 7285.7144
$ws.Range("I69").Value = 10166.667
$ws.Range("K69").Value = 30500.001
$ws.Range("M69").Value = -29626.001
$ws.Range("H72").Value = 7285.7144
$ws.Range("I72").Value = 10166.667
$ws.Range("K72").Value = 91500.003
$ws.Range("M72").Value = -87132.003
$ws.Range("H76").Value = 100002600
$ws.Range("I76").Value = 100002600
$ws.Range("K76").Value = 100002600
$ws.Range("M76").Value = -100002285
$ws.Range("H79").Value = 100002600
$ws.Range("I79").Value = 100002600
$ws.Range("K79").Value = 100002600
$ws.Range("M79").Value = -100001508
$ws.Range("H129").Value = 1255002.1
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1255002.1
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3765006.3
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -3775006.3
$ws.Range("H138").Value = 5245.0356
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 5245.0356
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 15735.1068
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -26015.1068

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5575425.5
$ws.Range("I61").Value = 2926526.5
$ws.Range("J61").Value = 29415514
$ws.Range("K61").Value = 2926526.5
$ws.Range("L61").Value = 29415514
$ws.Range("M61").Value = -2926314.5
$ws.Range("N61").Value = -29415938
$ws.Range("H63").Value = 1976.1904
$ws.Range("I63").Value = 1925
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1925
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1239
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 1976.1904
$ws.Range("I66").Value = 1925
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 9625
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -6193
$ws.Range("N66").Value = -21864
$ws.Range("H88").Value = 2989.7917
$ws.Range("I88").Value = 1340.6
$ws.Range("J88").Value = 4167.7856
$ws.Range("K88").Value = 1340.6
$ws.Range("L88").Value = 4167.7856
$ws.Range("M88").Value = -934.5999999999999
$ws.Range("N88").Value = -4979.7856
$ws.Range("H91").Value = 2989.7917
$ws.Range("I91").Value = 1340.6
$ws.Range("J91").Value = 4167.7856
$ws.Range("K91").Value = 1340.6
$ws.Range("L91").Value = 4167.7856
$ws.Range("M91").Value = 63.40000000000009
$ws.Range("N91").Value = -6975.7856
$ws.Range("H132").Value = 22488582
$ws.Range("I132").Value = 25547274
$ws.Range("J132").Value = 11366063
$ws.Range("K132").Value = 76641822
$ws.Range("L132").Value = 34098189
$ws.Range("M132").Value = -76639292
$ws.Range("N132").Value = -34103249
$ws.Range("H136").Value = 5575425.5
$ws.Range("I136").Value = 2926526.5
$ws.Range("J136").Value = 29415514
$ws.Range("K136").Value = 8779579.5
$ws.Range("L136").Value = 88246542
$ws.Range("M136").Value = -8777029.5
$ws.Range("N136").Value = -88251642
$ws.Range("H140").Value = 79987.3
$ws.Range("J140").Value = 79987.3
$ws.Range("L140").Value = 79987.3
$ws.Range("N140").Value = -90347.3

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1597.9333
$ws.Range("I105").Value = 1467.0869
$ws.Range("K105").Value = 1467.0869
$ws.Range("M105").Value = 279.9131
$ws.Range("H107").Value = 385341.3
$ws.Range("I107").Value = 476942.16
$ws.Range("J107").Value = 617.8
$ws.Range("K107").Value = 476942.16
$ws.Range("L107").Value = 617.8
$ws.Range("M107").Value = -475022.16
$ws.Range("N107").Value = -4457.8
$ws.Range("H132").Value = 31970
$ws.Range("J132").Value = 31970
$ws.Range("L132").Value = 31970
$ws.Range("N132").Value = -42090
$ws.Range("H134").Value = 11860354
$ws.Range("I134").Value = 15797708
$ws.Range("J134").Value = 48292.184
$ws.Range("K134").Value = 47393124
$ws.Range("L134").Value = 144876.552
$ws.Range("M134").Value = -47390589
$ws.Range("N134").Value = -149946.552
$ws.Range("H140").Value = 78400
$ws.Range("J140").Value = 78400
$ws.Range("L140").Value = 78400
$ws.Range("N140").Value = -88760

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 546.4375
$ws.Range("I107").Value = 169.28572
$ws.Range("J107").Value = 839.7778
$ws.Range("K107").Value = 169.28572
$ws.Range("L107").Value = 839.7778
$ws.Range("M107").Value = 1750.71428
$ws.Range("N107").Value = -4679.7778

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12525838
$ws.Range("I131").Value = 55666944
$ws.Range("J131").Value = 1001.3871
$ws.Range("K131").Value = 167000832
$ws.Range("L131").Value = 3004.1613
$ws.Range("M131").Value = -166995792
$ws.Range("N131").Value = -13084.1613
$ws.Range("H140").Value = 4473.5884
$ws.Range("I140").Value = 5034.9165
$ws.Range("K140").Value = 15104.7495
$ws.Range("M140").Value = -9924.749500000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6628528.5
$ws.Range("I70").Value = 2608397.5
$ws.Range("J70").Value = 20411834
$ws.Range("K70").Value = 2608397.5
$ws.Range("L70").Value = 20411834
$ws.Range("M70").Value = -2608127.5
$ws.Range("N70").Value = -20412374
$ws.Range("H73").Value = 6628528.5
$ws.Range("I73").Value = 2608397.5
$ws.Range("J73").Value = 20411834
$ws.Range("K73").Value = 2608397.5
$ws.Range("L73").Value = 20411834
$ws.Range("M73").Value = -2607461.5
$ws.Range("N73").Value = -20413706
$ws.Range("H80").Value = 11435.55
$ws.Range("I80").Value = 5960.5
$ws.Range("K80").Value = 5960.5
$ws.Range("M80").Value = -4962.5
$ws.Range("H83").Value = 11435.55
$ws.Range("I83").Value = 5960.5
$ws.Range("K83").Value = 29802.5
$ws.Range("M83").Value = -24810.5
$ws.Range("H102").Value = 2942.75
$ws.Range("I102").Value = 3381.3513
$ws.Range("J102").Value = 1860.8667
$ws.Range("K102").Value = 3381.3513
$ws.Range("L102").Value = 1860.8667
$ws.Range("M102").Value = -1759.3513
$ws.Range("N102").Value = -5104.8667
$ws.Range("H107").Value = 320.69232
$ws.Range("I107").Value = 117
$ws.Range("J107").Value = 558.3333
$ws.Range("K107").Value = 117
$ws.Range("L107").Value = 558.3333
$ws.Range("M107").Value = 1803
$ws.Range("N107").Value = -4398.3333
$ws.Range("H122").Value = 77787.734
$ws.Range("I122").Value = 113270.5
$ws.Range("J122").Value = 6822.2
$ws.Range("K122").Value = 339811.5
$ws.Range("L122").Value = 20466.6
$ws.Range("M122").Value = -337361.5
$ws.Range("N122").Value = -25366.6
$ws.Range("H132").Value = 13022611
$ws.Range("I132").Value = 15478279
$ws.Range("K132").Value = 46434837
$ws.Range("M132").Value = -46432307

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3705227.2
$ws.Range("I40").Value = 4274950.5
$ws.Range("J40").Value = 2025
$ws.Range("K40").Value = 4274950.5
$ws.Range("L40").Value = 2025
$ws.Range("M40").Value = -4274814.5
$ws.Range("N40").Value = -2297
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H135").Value = 50827.5
$ws.Range("J135").Value = 50827.5
$ws.Range("L135").Value = 50827.5
$ws.Range("N135").Value = -60967.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5634.59
$ws.Range("I107").Value = 7790.222
$ws.Range("J107").Value = 784.4167
$ws.Range("K107").Value = 23370.666
$ws.Range("L107").Value = 2353.2501
$ws.Range("M107").Value = -21450.666
$ws.Range("N107").Value = -6193.2501
$ws.Range("H109").Value = 42418
$ws.Range("J109").Value = 42418
$ws.Range("L109").Value = 42418
$ws.Range("N109").Value = -45192

